# Update Fonds de solidarite data with 2020-08-06 figures
# Updates columns C (nombre_aides) and D (montant_total) for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{Row=2; C=37570; D=54340619},
    @{Row=3; C=90658; D=132907979},
    @{Row=4; C=31086; D=46039022},
    @{Row=5; C=8665; D=12879203},
    @{Row=6; C=1984; D=2948506},
    @{Row=7; C=152; D=223093},
    @{Row=11; C=41125; D=55811165},
    @{Row=12; C=9616; D=13908599},
    @{Row=13; C=25870; D=37942243},
    @{Row=14; C=8298; D=12315218},
    @{Row=15; C=2142; D=3185383},
    @{Row=17; C=33; D=49500},
    @{Row=19; C=10185; D=13492137},
    @{Row=20; C=13340; D=19264971},
    @{Row=21; C=31565; D=46325192},
    @{Row=22; C=10203; D=15167555},
    @{Row=23; C=2628; D=3907182},
    @{Row=26; C=11640; D=15550098},
    @{Row=27; C=7617; D=11036435},
    @{Row=28; C=22411; D=32895431},
    @{Row=29; C=7790; D=11592633},
    @{Row=31; C=366; D=546415},
    @{Row=33; C=8266; D=10924364},
    @{Row=34; C=3222; D=4650137},
    @{Row=35; C=7793; D=11380820},
    @{Row=36; C=3167; D=4693461},
    @{Row=37; C=826; D=1230323},
    @{Row=38; C=162; D=241232},
    @{Row=40; C=2453; D=3314816},
    @{Row=41; C=17166; D=24828249},
    @{Row=42; C=50922; D=74654995},
    @{Row=43; C=18959; D=28162250},
    @{Row=44; C=5589; D=8323478},
    @{Row=45; C=1196; D=1784545},
    @{Row=49; C=16617; D=22133272},
    @{Row=50; C=1993; D=2891562},
    @{Row=51; C=6816; D=10020837},
    @{Row=52; C=2330; D=3479918},
    @{Row=56; C=6797; D=9361295},
    @{Row=57; C=926; D=1358579},
    @{Row=58; C=2307; D=3420317},
    @{Row=59; C=923; D=1374001},
    @{Row=60; C=316; D=473758},
    @{Row=61; C=99; D=148500},
    @{Row=63; C=1363; D=1916956},
    @{Row=64; C=15282; D=22077171},
    @{Row=65; C=44535; D=65175003},
    @{Row=66; C=15655; D=23267546},
    @{Row=67; C=4558; D=6789292},
    @{Row=68; C=915; D=1360668},
    @{Row=72; C=15031; D=19825086},
    @{Row=73; C=51014; D=74241233},
    @{Row=74; C=145141; D=213838917},
    @{Row=75; C=63313; D=94345270},
    @{Row=76; C=20213; D=30200209},
    @{Row=77; C=4778; D=7138723},
    @{Row=78; C=261; D=386670},
    @{Row=80; C=14; D=20125},
    @{Row=84; C=50451; D=68663346},
    @{Row=85; C=4558; D=6603884},
    @{Row=86; C=11506; D=16905652},
    @{Row=87; C=3864; D=5758906},
    @{Row=88; C=1341; D=2003989},
    @{Row=92; C=5363; D=7211575},
    @{Row=93; C=1585; D=2282432},
    @{Row=94; C=5117; D=7537139},
    @{Row=95; C=1933; D=2879437},
    @{Row=97; C=181; D=270613},
    @{Row=100; C=3509; D=4645963},
    @{Row=102; C=347; D=518030},
    @{Row=106; C=10710; D=15543478},
    @{Row=107; C=29100; D=42759163},
    @{Row=108; C=9756; D=14507650},
    @{Row=110; C=487; D=725546},
    @{Row=113; C=9743; D=12876321},
    @{Row=114; C=30279; D=43667790},
    @{Row=115; C=65903; D=96455711},
    @{Row=116; C=21296; D=31648832},
    @{Row=117; C=6040; D=8999826},
    @{Row=123; C=25718; D=34366715},
    @{Row=124; C=35793; D=51667188},
    @{Row=125; C=76476; D=111841010},
    @{Row=126; C=23751; D=35252389},
    @{Row=127; C=6366; D=9460051},
    @{Row=128; C=1225; D=1821911},
    @{Row=132; C=31611; D=41995351},
    @{Row=133; C=13166; D=19060394},
    @{Row=134; C=32202; D=47301440},
    @{Row=135; C=11443; D=17002542},
    @{Row=136; C=2948; D=4394805},
    @{Row=137; C=494; D=734990},
    @{Row=140; C=10771; D=14366304},
    @{Row=141; C=34858; D=50343506},
    @{Row=142; C=80895; D=118529191},
    @{Row=143; C=24260; D=36046706},
    @{Row=144; C=6362; D=9492567},
    @{Row=145; C=1427; D=2122730},
    @{Row=148; C=29030; D=39182221}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
